$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 8-10 (the MuSCs sending-cluster rows removed in the new data)
$ws.Range("A8:T10").EntireRow.Delete()

# Updated numeric values for rows 2-7 (columns G,H,I,J,M,N,O,P,Q,R,S,T)
$ws.Range("G2").Value = 0.01272866666666667
$ws.Range("H2").Value = 0.038186
$ws.Range("I2").Value = 0.04860514461513396
$ws.Range("J2").Value = 0.04860514461513396
$ws.Range("M2").Value = 3.241087666666667
$ws.Range("N2").Value = 9.723262999999999
$ws.Range("O2").Value = 0.02486257877280725
$ws.Range("P2").Value = 0.02486257877280725
$ws.Range("Q2").Value = 0.04125472454644444
$ws.Range("R2").Value = 0.371292520918
$ws.Range("S2").Value = 0.001208449236757456
$ws.Range("T2").Value = 0.001208449236757456

$ws.Range("G3").Value = 0.01272866666666667
$ws.Range("H3").Value = 0.038186
$ws.Range("I3").Value = 0.04860514461513396
$ws.Range("J3").Value = 0.04860514461513396
$ws.Range("O3").Value = 0.02096124117795788
$ws.Range("P3").Value = 0.02096124117795788
$ws.Range("Q3").Value = 0.03478119622466666
$ws.Range("R3").Value = 0.313030766022
$ws.Range("S3").Value = 0.001018824158767344
$ws.Range("T3").Value = 0.001018824158767344

$ws.Range("G4").Value = 0.01272866666666667
$ws.Range("H4").Value = 0.038186
$ws.Range("I4").Value = 0.04860514461513396
$ws.Range("J4").Value = 0.04860514461513396
$ws.Range("M4").Value = 124.3864796666667
$ws.Range("N4").Value = 373.159439
$ws.Range("O4").Value = 0.9541761800492348
$ws.Range("P4").Value = 0.9541761800492349
$ws.Range("Q4").Value = 1.583274037517111
$ws.Range("R4").Value = 14.249466337654
$ws.Range("S4").Value = 0.04637787121960916
$ws.Range("T4").Value = 0.04637787121960916

$ws.Range("I5").Value = 0.9513948553848661
$ws.Range("J5").Value = 0.9513948553848661
$ws.Range("M5").Value = 3.241087666666667
$ws.Range("N5").Value = 9.723262999999999
$ws.Range("O5").Value = 0.02486257877280725
$ws.Range("P5").Value = 0.02486257877280725
$ws.Range("Q5").Value = 0.8075180725125556
$ws.Range("R5").Value = 7.267662652613
$ws.Range("S5").Value = 0.02365412953604979
$ws.Range("T5").Value = 0.0236541295360498

$ws.Range("I6").Value = 0.9513948553848661
$ws.Range("J6").Value = 0.9513948553848661
$ws.Range("O6").Value = 0.02096124117795788
$ws.Range("P6").Value = 0.02096124117795788
$ws.Range("S6").Value = 0.01994241701919054
$ws.Range("T6").Value = 0.01994241701919054

$ws.Range("I7").Value = 0.9513948553848661
$ws.Range("J7").Value = 0.9513948553848661
$ws.Range("M7").Value = 124.3864796666667
$ws.Range("N7").Value = 373.159439
$ws.Range("O7").Value = 0.9541761800492348
$ws.Range("P7").Value = 0.9541761800492349
$ws.Range("Q7").Value = 30.99093287110989
$ws.Range("R7").Value = 278.918395839989
$ws.Range("S7").Value = 0.9077983088296256
$ws.Range("T7").Value = 0.9077983088296258
